$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 490
$ws.Range("I12").Value = 505
$ws.Range("K12").Value = 505
$ws.Range("M12").Value = -335
$ws.Range("H28").Value = 1610.2727
$ws.Range("I28").Value = 589.5
$ws.Range("J28").Value = 4332.3335
$ws.Range("K28").Value = 589.5
$ws.Range("L28").Value = 4332.3335
$ws.Range("M28").Value = -104.5
$ws.Range("N28").Value = -5302.3335
$ws.Range("H76").Value = 12720
$ws.Range("I76").Value = 4533.3335
$ws.Range("K76").Value = 4533.3335
$ws.Range("M76").Value = -4218.3335
$ws.Range("H79").Value = 12720
$ws.Range("I79").Value = 4533.3335
$ws.Range("K79").Value = 4533.3335
$ws.Range("M79").Value = -3441.3335
$ws.Range("H86").Value = 8036.3335
$ws.Range("I86").Value = 5975.2856
$ws.Range("K86").Value = 5975.2856
$ws.Range("M86").Value = -4852.2856
$ws.Range("H89").Value = 8036.3335
$ws.Range("I89").Value = 5975.2856
$ws.Range("K89").Value = 29876.428
$ws.Range("M89").Value = -24260.428
$ws.Range("H100").Value = 3067.0908
$ws.Range("I100").Value = 2059.625
$ws.Range("K100").Value = 2059.625
$ws.Range("M100").Value = -1518.625
$ws.Range("H135").Value = 3583
$ws.Range("I135").Value = 3583
$ws.Range("K135").Value = 32247
$ws.Range("M135").Value = -29712
$ws.Range("H137").Value = 10682.923
$ws.Range("I137").Value = 3853.7273
$ws.Range("K137").Value = 11561.1819
$ws.Range("M137").Value = -9011.1819
$ws.Range("H141").Value = 6444.2856
$ws.Range("I141").Value = 5851.75
$ws.Range("K141").Value = 17555.25
$ws.Range("M141").Value = -12375.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 35000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H45").Value = 3254.4333
$ws.Range("I45").Value = 3402.682
$ws.Range("J45").Value = 2846.75
$ws.Range("K45").Value = 3402.682
$ws.Range("L45").Value = 2846.75
$ws.Range("M45").Value = -3025.682
$ws.Range("N45").Value = -3600.75
$ws.Range("H110").Value = 4890.0713
$ws.Range("I110").Value = 6125.9
$ws.Range("J110").Value = 1800.5
$ws.Range("K110").Value = 6125.9
$ws.Range("L110").Value = 1800.5
$ws.Range("M110").Value = -4080.9
$ws.Range("N110").Value = -5890.5
$ws.Range("H132").Value = 14921.946
$ws.Range("I132").Value = 3867.9
$ws.Range("J132").Value = 27926.705
$ws.Range("K132").Value = 11603.7
$ws.Range("L132").Value = 83780.11500000001
$ws.Range("M132").Value = -9073.700000000001
$ws.Range("N132").Value = -88840.11500000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2005.7778
$ws.Range("I107").Value = 1925.4166
$ws.Range("K107").Value = 1925.4166
$ws.Range("M107").Value = -5.416600000000017
$ws.Range("H134").Value = 29054.4
$ws.Range("I134").Value = 31370.795
$ws.Range("J134").Value = 21894.637
$ws.Range("K134").Value = 94112.38499999999
$ws.Range("L134").Value = 65683.91099999999
$ws.Range("M134").Value = -91577.38499999999
$ws.Range("N134").Value = -70753.91099999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10790
$ws.Range("I86").Value = 13539.083
$ws.Range("J86").Value = 6666.375
$ws.Range("K86").Value = 13539.083
$ws.Range("L86").Value = 6666.375
$ws.Range("M86").Value = -12416.083
$ws.Range("N86").Value = -8912.375
$ws.Range("H89").Value = 10790
$ws.Range("I89").Value = 13539.083
$ws.Range("J89").Value = 6666.375
$ws.Range("K89").Value = 67695.41500000001
$ws.Range("L89").Value = 33331.875
$ws.Range("M89").Value = -62079.41500000001
$ws.Range("N89").Value = -44563.875
$ws.Range("H99").Value = 147428.11
$ws.Range("I99").Value = 220481.7
$ws.Range("J99").Value = 9438
$ws.Range("K99").Value = 220481.7
$ws.Range("L99").Value = 9438
$ws.Range("M99").Value = -218983.7
$ws.Range("N99").Value = -12434
$ws.Range("H126").Value = 147428.11
$ws.Range("I126").Value = 220481.7
$ws.Range("J126").Value = 9438
$ws.Range("K126").Value = 661445.1000000001
$ws.Range("L126").Value = 28314
$ws.Range("M126").Value = -658975.1000000001
$ws.Range("N126").Value = -33254

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 247.40909
$ws.Range("I2").Value = 337.41666
$ws.Range("J2").Value = 139.4
$ws.Range("K2").Value = 2024.49996
$ws.Range("L2").Value = 836.4000000000001
$ws.Range("M2").Value = -1911.49996
$ws.Range("N2").Value = -1062.4
$ws.Range("H23").Value = 637
$ws.Range("I23").Value = 700
$ws.Range("J23").Value = 611.8
$ws.Range("K23").Value = 2100
$ws.Range("L23").Value = 1835.4
$ws.Range("M23").Value = -1865
$ws.Range("N23").Value = -2305.4
$ws.Range("H64").Value = 6928.25
$ws.Range("J64").Value = 6928.25
$ws.Range("L64").Value = 20784.75
$ws.Range("N64").Value = -21324.75
$ws.Range("H67").Value = 6928.25
$ws.Range("J67").Value = 6928.25
$ws.Range("L67").Value = 20784.75
$ws.Range("N67").Value = -22656.75
$ws.Range("H103").Value = 1187.8889
$ws.Range("I103").Value = 172.8
$ws.Range("J103").Value = 2456.75
$ws.Range("K103").Value = 518.4000000000001
$ws.Range("L103").Value = 7370.25
$ws.Range("M103").Value = 360.5999999999999
$ws.Range("N103").Value = -9128.25
$ws.Range("H107").Value = 2095.697
$ws.Range("I107").Value = 561.8333
$ws.Range("J107").Value = 2972.1904
$ws.Range("K107").Value = 1685.4999
$ws.Range("L107").Value = 8916.5712
$ws.Range("M107").Value = 234.5001
$ws.Range("N107").Value = -12756.5712
$ws.Range("H137").Value = 4913.5
$ws.Range("J137").Value = 4496.5
$ws.Range("L137").Value = 13489.5
$ws.Range("N137").Value = -23689.5
$ws.Range("H140").Value = 2883.889
$ws.Range("I140").Value = 2744.5
$ws.Range("K140").Value = 8233.5
$ws.Range("M140").Value = -3053.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18333.8
$ws.Range("I70").Value = 17016.953
$ws.Range("K70").Value = 17016.953
$ws.Range("M70").Value = -16746.953
$ws.Range("H73").Value = 18333.8
$ws.Range("I73").Value = 17016.953
$ws.Range("K73").Value = 17016.953
$ws.Range("M73").Value = -16080.953
$ws.Range("H80").Value = 11628.037
$ws.Range("I80").Value = 9585.888999999999
$ws.Range("K80").Value = 9585.888999999999
$ws.Range("M80").Value = -8587.888999999999
$ws.Range("H83").Value = 11628.037
$ws.Range("I83").Value = 9585.888999999999
$ws.Range("K83").Value = 47929.44499999999
$ws.Range("M83").Value = -42937.44499999999
$ws.Range("H107").Value = 1836.6666
$ws.Range("I107").Value = 1988.5
$ws.Range("K107").Value = 1988.5
$ws.Range("M107").Value = -68.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2579.0571
$ws.Range("J22").Value = 3535.1428
$ws.Range("L22").Value = 3535.1428
$ws.Range("N22").Value = -4125.1428
$ws.Range("H27").Value = 2579.0571
$ws.Range("J27").Value = 3535.1428
$ws.Range("L27").Value = 3535.1428
$ws.Range("N27").Value = -3749.1428
$ws.Range("H46").Value = 3819.1
$ws.Range("I46").Value = 2200
$ws.Range("J46").Value = 3999
$ws.Range("K46").Value = 2200
$ws.Range("L46").Value = 3999
$ws.Range("M46").Value = -2012
$ws.Range("N46").Value = -4375
$ws.Range("H61").Value = 4061.2856
$ws.Range("I61").Value = 3036.476
$ws.Range("J61").Value = 7135.7144
$ws.Range("K61").Value = 3036.476
$ws.Range("L61").Value = 7135.7144
$ws.Range("M61").Value = -2834.476
$ws.Range("N61").Value = -7539.7144
$ws.Range("H93").Value = 5495.0645
$ws.Range("I93").Value = 4210.25
$ws.Range("J93").Value = 7831.091
$ws.Range("K93").Value = 4210.25
$ws.Range("L93").Value = 7831.091
$ws.Range("M93").Value = -2962.25
$ws.Range("N93").Value = -10327.091
$ws.Range("H113").Value = 4061.2856
$ws.Range("I113").Value = 3036.476
$ws.Range("J113").Value = 7135.7144
$ws.Range("K113").Value = 3036.476
$ws.Range("L113").Value = 7135.7144
$ws.Range("M113").Value = -866.4760000000001
$ws.Range("N113").Value = -11475.7144
$ws.Range("H122").Value = 6308.763
$ws.Range("J122").Value = 6667.579
$ws.Range("L122").Value = 20002.737
$ws.Range("N122").Value = -24902.737
$ws.Range("H130").Value = 140888
$ws.Range("J130").Value = 140888
$ws.Range("L130").Value = 140888
$ws.Range("N130").Value = -150928
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 125000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 125000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 125000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -135280

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 40285.57
$ws.Range("I62").Value = 50000
$ws.Range("J62").Value = 36399.8
$ws.Range("K62").Value = 50000
$ws.Range("L62").Value = 36399.8
$ws.Range("M62").Value = -49376
$ws.Range("N62").Value = -37647.8
$ws.Range("H65").Value = 40285.57
$ws.Range("I65").Value = 50000
$ws.Range("J65").Value = 36399.8
$ws.Range("K65").Value = 250000
$ws.Range("L65").Value = 181999
$ws.Range("M65").Value = -246880
$ws.Range("N65").Value = -188239
$ws.Range("H81").Value = 38999.668
$ws.Range("I81").Value = 38999.668
$ws.Range("K81").Value = 77999.336
$ws.Range("M81").Value = -76938.336
$ws.Range("H84").Value = 38999.668
$ws.Range("I84").Value = 38999.668
$ws.Range("K84").Value = 389996.68
$ws.Range("M84").Value = -384692.68
$ws.Range("H136").Value = 11911.156
$ws.Range("I136").Value = 1323.9
$ws.Range("J136").Value = 29556.584
$ws.Range("K136").Value = 3971.7
$ws.Range("L136").Value = 88669.75199999999
$ws.Range("M136").Value = -1421.7
$ws.Range("N136").Value = -93769.75199999999
